# Refresh scraped crypto price/symbol data (GitHub Actions update, 2022-12-26)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'243.07"

# Row 3
$ws.Range("D3").Value = "'23.00"

# Row 4
$ws.Range("D4").Value = "'5.412"

# Row 5
$ws.Range("D5").Value = "'0.05921"

# Row 6
$ws.Range("D6").Value = "'3.439"

# Row 7
$ws.Range("D7").Value = "'6.515"

# Row 8
$ws.Range("D8").Value = "'0.8096"

# Row 9
$ws.Range("D9").Value = "'0.9284"

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1431"
$ws.Range("E10").Value = '9WazirXWRX'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.07412"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

# Row 12
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = "'0.03246"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.03080"
$ws.Range("E13").Value = '12BitrueCoinBTR'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.09360"
$ws.Range("E14").Value = '13BitMartTokenBMX'

# Row 15
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = "'3.856"
$ws.Range("E15").Value = '14MCDexMCB'

# Row 16
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = "'0.001580"
$ws.Range("E16").Value = '15BitForexTokenBF'

# Row 17
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = "'0.04676"
$ws.Range("E17").Value = '16CoinExTokenCET'

# Row 18
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = "'0.0005892"
$ws.Range("E18").Value = '17OneONEWorstin24h'

# Row 19
$ws.Range("D19").Value = "'0.006006"

# Row 21
$ws.Range("D21").Value = "'0.004911"

# Row 22
$ws.Range("D22").Value = "'0.00006803"

# Row 24
$ws.Range("D24").Value = "'2.124"

# Row 26
$ws.Range("D26").Value = "'0.1303"

# Row 27
$ws.Range("D27").Value = "'0.0002302"

# Row 40
$ws.Range("D40").Value = "'0.03961"

# Row 41
$ws.Range("D41").Value = "'0.006433"

# Row 43
$ws.Range("D43").Value = "'0.002571"

# Row 44
$ws.Range("D44").Value = "'0.008756"

# Row 45
$ws.Range("D45").Value = "'0.00005232"

# Row 46
$ws.Range("D46").Value = "'0.00000000750"

# Row 47
$ws.Range("D47").Value = "'0.6702"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'

# Row 48
$ws.Range("D48").Value = "'0.002395"

# Row 49
$ws.Range("D49").Value = "'0.00002101"

# Row 50
$ws.Range("D50").Value = "'0.0002000"
